$wb = $excel.ActiveWorkbook

# --- Replace status text "Ready for handoff" -> "In Translation" everywhere it occurs ---
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $found = $used.Find("Ready for handoff")
    if ($found -ne $null) {
        $firstAddress = $found.Address()
        do {
            $found.Value = "In Translation"
            $found = $used.FindNext($found)
        } while ($found -ne $null -and $found.Address() -ne $firstAddress)
    }
}

# --- Narrow the Status column(s) to their new width ---
# Original XML width 17.2159881591797 -> new XML width 13.4101845877511.
# The closest value reachable through the ColumnWidth COM property (which
# snaps to 1/6-character increments) is 13.333333333333334, produced by
# assigning 12.5.
$newColumnWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth  # column E
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth  # column F

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C
